# Mathison.docx edit:
#   "Functions" (underlined heading, paragraph 37) -> "Function Commands"
#   split across two runs ("Function" / " Commands"), and the document's
#   singleton "_GoBack" bookmark (previously sitting after the "END" /
#   page-break run near the end of the Functions section) moves to sit
#   right after the newly typed " Commands" text - mirroring what Word
#   does automatically when you edit text at a new location.

$d = $word.ActiveDocument

# Locate the paragraph that holds the "Functions" heading.
$targetPar = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -eq "Functions`r") {
        $targetPar = $d.Paragraphs($i)
        break
    }
}
if ($targetPar -eq $null) {
    throw "Could not locate the 'Functions' heading paragraph"
}

$parStart = $targetPar.Range.Start
# "Functions" text runs from $parStart to $parStart+9 (9 letters); the
# trailing "s" is the 9th character.
$sStart = $parStart + 8
$sEnd = $parStart + 9

# 1) Delete the trailing "s" so the run reads "Function".
$sRange = $d.Range($sStart, $sEnd)
$sRange.Text = ""

# 2) Type " Commands" right after "Function". Word would normally coalesce
#    this straight back into the same run as "Function" because the
#    formatting is identical, so we immediately wall the two apart with a
#    temporary bookmark (placed mid-run, a safe, non-boundary position) and
#    then discard the bookmark - the run split survives even once the
#    bookmark is removed.
$insertPoint = $d.Range($sStart, $sStart)
$insertPoint.InsertAfter(" Commands")

$d.Bookmarks.Add("_GoBack", $d.Range($sStart, $sStart))
$d.Bookmarks("_GoBack").Delete()

# 3) Re-create "_GoBack" at the true end of the paragraph's text (right
#    after " Commands", before the paragraph mark). Placing a zero-length
#    bookmark exactly at that boundary needs a one-character placeholder
#    scratch run to land on, which is then deleted, leaving a clean
#    zero-width bookmark in place.
$endPos = $sStart + " Commands".Length
$placeholder = $d.Range($endPos, $endPos)
$placeholder.InsertAfter("X")
$d.Bookmarks.Add("_GoBack", $d.Range($endPos, $endPos))
$d.Range($endPos, $endPos + 1).Text = ""
